$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix separator punctuation in "Razon social" / "Nombre Fantasia" names ---
$nameEdits = @(
    @('E65', 'PITTER ROLANDO LJ. CERGNEUX MARIO M Y CERGNEUX DANIEL F  SH'),
    @('F65', 'PITTER ROLANDO LJ. CERGNEUX MARIO M Y CERGNEUX DANIEL F  SH'),
    @('E69', 'FERNANDEZ MARIO H. GALLICET OSCAR M'),
    @('E70', 'IZAGUIRRE CARLOS MARIA. MOREND MARIA ELENA Y MOREND MARIA TERESA'),
    @('F70', 'IZAGUIRRE CARLOS MARIA. MOREND MARIA ELENA Y MOREND MARIA TERESA'),
    @('E72', 'MARSICO GUILLERMO MIGUEL. MARSICO JUAN EDUARDO'),
    @('E90', 'FERNANDEZ MARIO H. GALLICET OSCAR M'),
    @('E91', 'IZAGUIRRE CARLOS MARIA. MOREND MARIA ELENA Y MOREND MARIA TERESA'),
    @('F91', 'IZAGUIRRE CARLOS MARIA. MOREND MARIA ELENA Y MOREND MARIA TERESA'),
    @('E103', 'RICCOTTI. MARIANA EDITH')
)
foreach ($edit in $nameEdits) {
    $ws.Range($edit[0]).Value = $edit[1]
}

# --- Reformat "Importe" amounts from "1.234,56" (Argentine) style to plain "1234.56" ---
# These cells hold the amount as TEXT (shared string), so we must force a Text number
# format before assignment; otherwise Excel auto-parses the digit string into a real
# number and drops the trailing zero(s)/formatting.
$amountEdits = @(
    @('H2', '6730.00'),
    @('H3', '4016.00'),
    @('H4', '80.00'),
    @('H5', '50001.00'),
    @('H6', '160.00'),
    @('H7', '3896.68'),
    @('H8', '1170.97'),
    @('H9', '2559.86'),
    @('H10', '85539.87'),
    @('H11', '82130.33'),
    @('H12', '13124.00'),
    @('H13', '841.98'),
    @('H14', '11614.75'),
    @('H15', '1240.50'),
    @('H16', '3874.77'),
    @('H17', '1622.75'),
    @('H18', '4956.75'),
    @('H19', '6342.83'),
    @('H20', '211.20'),
    @('H21', '1443.00'),
    @('H22', '128.73'),
    @('H23', '6031.65'),
    @('H24', '200.00'),
    @('H25', '950.00'),
    @('H26', '6200.10'),
    @('H27', '411.21'),
    @('H28', '283.52'),
    @('H29', '9.00'),
    @('H30', '817.76'),
    @('H31', '1694.60'),
    @('H32', '4831.20'),
    @('H33', '1168.80'),
    @('H34', '135.00'),
    @('H35', '1274.00'),
    @('H36', '308.77'),
    @('H37', '3629.00'),
    @('H38', '7331.85'),
    @('H39', '630.00'),
    @('H40', '1727.36'),
    @('H41', '141.00'),
    @('H42', '2454.87'),
    @('H43', '505.55'),
    @('H44', '704.20'),
    @('H45', '12910.40'),
    @('H46', '150.00'),
    @('H47', '168.50'),
    @('H48', '1.32'),
    @('H49', '6.50'),
    @('H50', '19.36'),
    @('H51', '162.33'),
    @('H52', '147.96'),
    @('H53', '986.24'),
    @('H54', '200.00'),
    @('H55', '1298.20'),
    @('H56', '3030.00'),
    @('H57', '7838.00'),
    @('H58', '92.00'),
    @('H59', '1470.00'),
    @('H60', '4889.20'),
    @('H61', '800.00'),
    @('H62', '4602.00'),
    @('H63', '1390.00'),
    @('H64', '5700.00'),
    @('H65', '735.00'),
    @('H66', '300.00'),
    @('H67', '10853.20'),
    @('H68', '1747.30'),
    @('H69', '373.00'),
    @('H70', '781.23'),
    @('H71', '1164.12'),
    @('H72', '750.00'),
    @('H73', '46.00'),
    @('H74', '340.00'),
    @('H75', '2120.00'),
    @('H76', '6686.01'),
    @('H77', '10945.00'),
    @('H78', '3025.00'),
    @('H79', '6071.40'),
    @('H80', '64.00'),
    @('H81', '52.50'),
    @('H82', '409.19'),
    @('H83', '3213.00'),
    @('H84', '72.00'),
    @('H85', '1809.50'),
    @('H86', '1690.00'),
    @('H87', '310.00'),
    @('H88', '3500.00'),
    @('H89', '1280.00'),
    @('H90', '2631.60'),
    @('H91', '89.39'),
    @('H92', '722.00'),
    @('H93', '2103.50'),
    @('H94', '122.50'),
    @('H95', '6312.75'),
    @('H96', '808.00'),
    @('H97', '445.25'),
    @('H98', '133.19'),
    @('H99', '2977.00'),
    @('H100', '1387.34'),
    @('H101', '27597.00'),
    @('H102', '30.00'),
    @('H103', '12000.00'),
    @('H104', '180.00'),
    @('H105', '120.00'),
    @('H106', '3287.00'),
    @('H107', '258.00'),
    @('H108', '9580.00'),
    @('H109', '5270.00'),
    @('H110', '17590.00'),
    @('H111', '4235.00'),
    @('H112', '1300.00'),
    @('H113', '3700.00'),
    @('H114', '2070.00'),
    @('H115', '256918.50'),
    @('H116', '4500.00'),
    @('H117', '6282.00'),
    @('H118', '70.00'),
    @('H119', '3952.17'),
    @('H120', '249901.60'),
    @('H121', '240.00'),
    @('H122', '1381.00'),
    @('H123', '300.00'),
    @('H124', '900.00'),
    @('H125', '285.00'),
    @('H126', '642.00'),
    @('H127', '650.00'),
    @('H128', '300.00'),
    @('H129', '200.00'),
    @('H130', '1500.00'),
    @('H131', '1682.00'),
    @('H132', '6643.26'),
    @('H133', '800.00'),
    @('H134', '200.00'),
    @('H135', '220.00'),
    @('H136', '2173.00'),
    @('H137', '1700.00'),
    @('H138', '200.00'),
    @('H139', '62.50'),
    @('H140', '1200.09'),
    @('H141', '170.00'),
    @('H142', '139.00'),
    @('H143', '1900.00'),
    @('H144', '150.00'),
    @('H145', '60.00'),
    @('H146', '732.00'),
    @('H147', '185.00'),
    @('H148', '738.39'),
    @('H149', '2280.00'),
    @('H150', '967.26'),
    @('H151', '4796.00'),
    @('H152', '3050.00'),
    @('H153', '45.00'),
    @('H154', '215.00'),
    @('H155', '59.80'),
    @('H156', '2012.00'),
    @('H157', '4307.55'),
    @('H158', '12.00'),
    @('H159', '798.00'),
    @('H160', '802.00'),
    @('H161', '155.50'),
    @('H162', '380.00'),
    @('H163', '48.24'),
    @('H164', '5418.01'),
    @('H165', '884.44'),
    @('H166', '19.50'),
    @('H167', '898.60'),
    @('H168', '1762.25'),
    @('H169', '1404.00'),
    @('H170', '3259.74'),
    @('H171', '1926.62'),
    @('H172', '7380.00'),
    @('H173', '9900.00'),
    @('H174', '12879.00'),
    @('H175', '230.00'),
    @('H176', '293309.93'),
    @('H177', '9350.00'),
    @('H178', '1900.00'),
    @('H179', '1400.00'),
    @('H180', '3300.00'),
    @('H181', '4747.99'),
    @('H182', '1630.00')
)

$amountRange = $ws.Range("H2:H182")
$amountRange.NumberFormat = "@"
foreach ($edit in $amountEdits) {
    $ws.Range($edit[0]).Value = $edit[1]
}
$amountRange.Style = "Normal"

